$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.230988666666666
$ws.Range("H2").Value = 15.692966
$ws.Range("I2").Value = 0.2129406655351238
$ws.Range("J2").Value = 0.2129406655351238
$ws.Range("Q2").Value = 2.860768417261777
$ws.Range("R2").Value = 25.746915755356
$ws.Range("S2").Value = 0.2129406655351238
$ws.Range("T2").Value = 0.2129406655351238

# Row 3
$ws.Range("I3").Value = 0.345577477529236
$ws.Range("J3").Value = 0.3455774775292359
$ws.Range("S3").Value = 0.345577477529236
$ws.Range("T3").Value = 0.3455774775292359

# Row 4
$ws.Range("G4").Value = 8.418577333333333
$ws.Range("H4").Value = 25.255732
$ws.Range("I4").Value = 0.3426995496362334
$ws.Range("J4").Value = 0.3426995496362334
$ws.Range("Q4").Value = 4.604024533056889
$ws.Range("R4").Value = 41.436220797512
$ws.Range("S4").Value = 0.3426995496362334
$ws.Range("T4").Value = 0.3426995496362334

# Row 5
$ws.Range("G5").Value = 2.426634333333333
$ws.Range("H5").Value = 7.279902999999999
$ws.Range("I5").Value = 0.09878230729940689
$ws.Range("J5").Value = 0.09878230729940687
$ws.Range("Q5").Value = 1.327098815044222
$ws.Range("R5").Value = 11.943889335398
$ws.Range("S5").Value = 0.09878230729940689
$ws.Range("T5").Value = 0.09878230729940687
